$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.91%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'32.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.50%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-2.89%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.07914"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-2.56%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'2.105"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-13.81%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'7.872"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.56%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'3.810"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.74%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.9253"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.01%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1751"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.54%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07927"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'7.58%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.08714"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.82%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.03136"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.75%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.1004"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.29%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.001511"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.30%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.005967"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.01%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.467"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-3.57%"
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.276"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.39%"
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'0.74%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.1292"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-3.60%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'4.212"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.83%"
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'6.62%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04592"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.91%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001238"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-1.03%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.004467"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-1.61%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0001250"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'3.93%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D39").Value = "'0.01721"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-1.36%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.04806"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'4.53%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.007419"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'7.50%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.1363"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-1.14%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.002361"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'10.03%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.01026"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'4.03%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00006023"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-3.10%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.23%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.003393"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-59.71%"
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'2.31%"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.23%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.0002001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.23%"
$ws.Range("E50").Style = "Normal"
